$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the B, D, F values between row 8 and row 9 (data correction:
# the two entries had been transposed).
$b8 = $ws.Range("B8").Value()
$d8 = $ws.Range("D8").Value()
$f8 = $ws.Range("F8").Value()

$b9 = $ws.Range("B9").Value()
$d9 = $ws.Range("D9").Value()
$f9 = $ws.Range("F9").Value()

$ws.Range("B8").Value = $b9
$ws.Range("D8").Value = $d9
$ws.Range("F8").Value = $f9

$ws.Range("B9").Value = $b8
$ws.Range("D9").Value = $d8
$ws.Range("F9").Value = $f8

# Move the active selection to D13 (matches author's final cursor position).
$ws.Range("D13").Select()
